$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.0129867619495103
$ws.Range("B1").Value = -0.00477222198665538
$ws.Range("C1").Value = -0.00614888577083178
$ws.Range("D1").Value = 2.00109081513898
